$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10993.53
$ws.Range("M2").Value = 7329.02

$ws.Range("G3").Value = 186811.2
$ws.Range("M3").Value = 124540.8

$ws.Range("F4").Value = 87163.98
$ws.Range("G4").Value = 8788.620000000001
$ws.Range("M4").Value = 8788.620000000001

$ws.Range("G5").Value = 68624.48
$ws.Range("L5").Value = 22858.2
$ws.Range("M5").Value = 22814.55

$ws.Range("G6").Value = 41412.83
$ws.Range("L6").Value = 13825.33
$ws.Range("M6").Value = 13744.13

$ws.Range("F7").Value = 84311.00999999999

$ws.Range("F8").Value = 92527.36

$ws.Range("F9").Value = 91990.97

$ws.Range("F11").Value = 44189.58
$ws.Range("G11").Value = 42502.95
$ws.Range("L11").Value = 8500.59
$ws.Range("M11").Value = 8500.59

$ws.Range("F12").Value = 48312.09
$ws.Range("G12").Value = 42531.7
$ws.Range("L12").Value = 8506.34
$ws.Range("M12").Value = 8506.34

$ws.Range("G13").Value = 55183.14999999999
$ws.Range("M13").Value = 11036.63

$ws.Range("G14").Value = 55656.1
$ws.Range("M14").Value = 11131.22

$ws.Range("F15").Value = 4078500.49
$ws.Range("G15").Value = 960293.85
$ws.Range("M15").Value = 310028.22

$ws.Range("G16").Value = 5503.69
$ws.Range("M16").Value = 5503.69
